$d = $word.ActiveDocument
$t = $d.Tables.Item(4)
$cell = $t.Cell(3, 3)
$cell.Range.Text = "{{project_end_date}}"
